# Regenerate the "K" column (column G) values for the save_data sheet.
# The upstream data-generation script was changed to use the option "K" (strike)
# column values instead of the previous "Strike#" derivation, after
# recalculating std/mean and writing fresh simulated values (s_vals).
# This updates each data row's K (column G) with its newly computed value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (column G), for data rows 2..65.
# Row 58 is intentionally absent: its recomputed value is unchanged (1).
$kValues = [ordered]@{
    2  = 2
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 2
    14 = 1
    15 = 1
    16 = 0
    17 = 2
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 2
    23 = 2
    24 = 1
    25 = 1
    26 = 0
    27 = 1
    28 = 0
    29 = 0
    30 = 2
    31 = 1
    32 = 0
    33 = 0
    34 = 1
    35 = 3
    36 = 0
    37 = 1
    38 = 0
    39 = 2
    40 = 1
    41 = 1
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 0
    47 = 0
    48 = 1
    49 = 1
    50 = 2
    51 = 1
    52 = 0
    53 = 2
    54 = 1
    55 = 2
    56 = 1
    57 = 0
    59 = 1
    60 = 1
    61 = 3
    62 = 2
    63 = 2
    64 = 2
    65 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
